# Acta de Solicitud de Diferencia Laboral ("Generar acta de solicitud
# completo"): turn the two hand-typed placeholder runs into the full set
# of ${...} merge fields used by the template engine.
#
# Each replacement is applied as a raw WordprocessingML fragment via
# Range.InsertXML so the resulting paragraph keeps one run per merge
# field (mirroring how the template was authored), instead of Word
# collapsing everything into a single run the way Find/Replace would.

$d = $word.ActiveDocument

# Locate the two paragraphs that need rewriting by their current text,
# rather than assuming fixed indices.
$licdaParaIndex = $null
$bigParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($licdaParaIndex -eq $null -and $t -like "*LICDA*") {
        $licdaParaIndex = $i
    }
    if ($bigParaIndex -eq $null -and $t -like "*DIRECCIÓN GENERAL DE TRABAJO*") {
        $bigParaIndex = $i
    }
}

# --- "LICDA. ${delegado}" -> "${nombre_delegado}" -----------------------
# This also removes the _GoBack bookmark that used to sit inside this
# run; it is re-created below, inside the big paragraph, at the point
# where editing actually finished.
$licdaPara = $d.Paragraphs.Item($licdaParaIndex).Range
$licdaPara.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>${</w:t></w:r><w:r><w:t>nombre_</w:t></w:r><w:r><w:t>delegado}</w:t></w:r></w:p>')

# --- Big narrative paragraph: fill in every ${...} merge field -----------
# Replace the single placeholder-laden run with the fully tokenised
# text, split into one run per ${...} field (and per literal chunk of
# surrounding prose), and put the _GoBack bookmark back around
# ${mes_audiencia2}, matching where the edit actually finished.
$bigPara = $d.Paragraphs.Item($bigParaIndex).Range
$bigPara.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">EN LA DIRECCIÓN GENERAL DE TRABAJO: San Salvador, a las </w:t></w:r><w:r><w:t>${hora_expediente}</w:t></w:r><w:r><w:t xml:space="preserve"> horas y </w:t></w:r><w:r><w:t>${minuto_expediente}</w:t></w:r><w:r><w:t xml:space="preserve"> minutos del día </w:t></w:r><w:r><w:t>${dia_expediente}</w:t></w:r><w:r><w:t xml:space="preserve"> de </w:t></w:r><w:r><w:t>${mes_expediente}</w:t></w:r><w:r><w:t xml:space="preserve"> del </w:t></w:r><w:r><w:t>${anio_expediente}</w:t></w:r><w:r><w:t xml:space="preserve">. Comparecen los directivos sindicales señores: </w:t></w:r><w:r><w:t>${directivos}</w:t></w:r><w:r><w:t xml:space="preserve"> calidades que acreditan por medio de Credenciales, extendidas por el Departamento Nacional de Organizaciones Sociales,  de la  Dirección General de Trabajo, de este Ministerio, Salvadoreños, señalando para oír notificaciones en: </w:t></w:r><w:r><w:t>${direccion_sindicato}</w:t></w:r><w:r><w:t xml:space="preserve">; y DICEN: que solicitan la intervención de esta Dirección General para que se les nombre un Delegado y se cite en legal forma a </w:t></w:r><w:r><w:t>${nombre_</w:t></w:r><w:r><w:t>empresa</w:t></w:r><w:r><w:t>}</w:t></w:r><w:r><w:t xml:space="preserve">, representada legalmente por </w:t></w:r><w:r><w:t>${representante_legal}</w:t></w:r><w:r><w:t xml:space="preserve">, a quien se le puede notificar y citar en: </w:t></w:r><w:r><w:t>${direccion_empresa}</w:t></w:r><w:r><w:t xml:space="preserve">. Para que en audiencia conciliatoria que se llevará a cabo en estas oficinas, ubicadas en: ALAMEDA JUAN PABLO II, Y DIECISIETE AVENIDA NORTE, EDIFICIO DOS, PRIMER NIVEL,  EX OFICINAS SETEFE, PLAN MAESTRO CENTRO DE GOBIERNO, DE ESTA CIUDAD, se procure resolver la presente diferencia laboral, la cual consiste en:  </w:t></w:r><w:r><w:t>${motivo}</w:t></w:r><w:r><w:t xml:space="preserve">. La suscrita Directora General de Trabajo RESUELVE: Admitir la presente solicitud de conformidad a lo establecido en el Artículo veinticuatro de la Ley de Organización y Funciones del Sector Trabajo y Previsión Social; tener por parte en las presentes diligencias a los directivos Sindicales antes mencionados, y designase como Delegados de la suscrita para intervenir en tales diligencias a los Licenciados </w:t></w:r><w:r><w:t>${nombre_delegado}</w:t></w:r><w:r><w:t xml:space="preserve">, CÍTESE POR PRIMERA VEZ a </w:t></w:r><w:r><w:t>${nombre_empresa}</w:t></w:r><w:r><w:t xml:space="preserve">, por medio de su Representante Legal </w:t></w:r><w:r><w:t>${representante_legal}</w:t></w:r><w:r><w:t xml:space="preserve">, para que comparezca a este Dirección General de Trabajo, a las  </w:t></w:r><w:r><w:t>${hora_audiencia}</w:t></w:r><w:r><w:t xml:space="preserve"> horas y </w:t></w:r><w:r><w:t>${minuto_audiencia}</w:t></w:r><w:r><w:t xml:space="preserve"> minutos del día </w:t></w:r><w:r><w:t>${dia_audiencia}</w:t></w:r><w:r><w:t xml:space="preserve"> del mes de </w:t></w:r><w:r><w:t>${mes_audiencia}</w:t></w:r><w:r><w:t xml:space="preserve"> del corriente año, para celebrar audiencia conciliatoria con los solicitantes. De no verificarse la audiencia conciliatoria en la primera cita señalada, debido a la inasistencia del patrono, no obstante haber sido notificado y citado legalmente. CÍTESELE POR SEGUNDA VEZ, para que comparezca a estas oficinas a las </w:t></w:r><w:r><w:t>${hora_audiencia2}</w:t></w:r><w:r><w:t xml:space="preserve"> horas y </w:t></w:r><w:r><w:t>${minuto_audiencia2}</w:t></w:r><w:r><w:t xml:space="preserve"> minutos del día </w:t></w:r><w:r><w:t>${dia_audiencia</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t>}</w:t></w:r><w:r><w:t xml:space="preserve"> del mes de </w:t></w:r><w:r><w:t>${mes_audiencia</w:t></w:r><w:r><w:t>2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>}</w:t></w:r><w:r><w:t xml:space="preserve"> del corriente año a fin de llevar a cabo la audiencia conciliatoria antes mencionada; PREVINIÉNDOLE que de no comparecer a este segundo señalamiento, incurrirá en la MULTA que señala el Artículo treinta y dos de la Ley de Organización y Funciones del Sector Trabajo y Previsión Social, si el solicitado es una persona jurídica, deberá comparecer el representante legal con documentación que acredite fecha recientemente su personería y en caso de no poder comparecer personalmente deberá hacerlo por medio de apoderado de conformidad con lo establecido en el Artículo trescientos setenta y cinco del Código de Trabajo vigente. Los comparecientes se da por notificados y citados de los señalamientos anteriores. No habiendo nada más que hacer constar, se da por terminada la presente acta y leída que les fue a los comparecientes la ratifican y para constancia firmamos.</w:t></w:r></w:p>')
